$wb = $excel.ActiveWorkbook

# --- Sheet 1: insert new row 2 and shift existing rows down ---
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "日期：2021/12/09"
$ws.Range("B2").Value = "202201"
$ws.Range("C2").Value = 17852
$ws.Range("D2").Value = 12187
$ws.Range("E2").Value = 18084076
$ws.Range("F2").Value = 17664

# --- Sheet 2: insert new row 2 and shift existing rows down ---
$ws = $wb.Worksheets.Item(2)
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "日期：2021/12/09"
$ws.Range("B2").Value = -0.03

# --- Sheet 3: insert new row 2 and shift existing rows down ---
$ws = $wb.Worksheets.Item(3)
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "110年12月09日"
$ws.Range("B2").Value = 61.47
$ws.Range("C2").Value = 20.23

# --- Sheet 4: insert new row 2 and shift existing rows down ---
$ws = $wb.Worksheets.Item(4)
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "110年12月09日"
$ws.Range("B2").Value = 17876.02

# --- Sheet 5: insert new row 2 and shift existing rows down ---
$ws = $wb.Worksheets.Item(5)
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "2021/12/09"
$ws.Range("B2").Value = 47451
$ws.Range("C2").Value = 56262
$ws.Range("D2").Value = -545
$ws.Range("E2").Value = 789
$ws.Range("F2").Value = 25804
$ws.Range("G2").Value = 50404
$ws.Range("H2").Value = 803
$ws.Range("I2").Value = -11
$ws.Range("J2").Value = -24600
$ws.Range("K2").Value = 814
$ws.Range("L2").Value = -1348
$ws.Range("M2").Value = 800
$ws.Range("N2").Value = -2148
